$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-7 from 2023-11-03 (45233)
# to 2023-11-13 (45243).
$ws.Range("C2:C7").Value = 45243
